$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$CellAddress,
        [string]$NewValue
    )
    $rng = $ws.Range($CellAddress)
    $rng.NumberFormat = "@"
    $rng.Value = $NewValue
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '36.805.26'
Set-TextValue 'D3' '1.969.64'
Set-TextValue 'E3' '  +1.18%  '
Set-TextValue 'E4' '  -0.15%  '
Set-TextValue 'D5' '245.14'
Set-TextValue 'E5' '  +0.47%  '
Set-TextValue 'D6' '0.621'
Set-TextValue 'E6' '  +0.84%  '
Set-TextValue 'D7' '59.14'
Set-TextValue 'E7' '  +0.93%  '
Set-TextValue 'E8' '  -0.07%  '
Set-TextValue 'D9' '0.375'
Set-TextValue 'E9' '  +1.69%  '
Set-TextValue 'D10' '0.0813'
Set-TextValue 'E10' '  -3.22%  '
Set-TextValue 'E11' '  -0.19%  '
Set-TextValue 'D12' '22.56'
Set-TextValue 'E12' '  +3.84%  '
Set-TextValue 'D13' '2.255.70'
Set-TextValue 'E13' '  +1.13%  '
Set-TextValue 'D14' '0.829'
Set-TextValue 'E14' '  +0.03%  '
Set-TextValue 'D15' '13.78'
Set-TextValue 'E15' '  +1.11%  '
Set-TextValue 'E16' '  +0.20%  '
Set-TextValue 'D17' '1.963.38'
Set-TextValue 'E17' '  +0.82%  '
Set-TextValue 'D18' '36.716.38'
Set-TextValue 'E18' '  +0.75%  '
Set-TextValue 'D19' '69.89'
Set-TextValue 'E19' '  +0.12%  '
Set-TextValue 'D20' '0.0₃0863'
Set-TextValue 'E20' '  -1.07%  '
Set-TextValue 'E21' '  +1.55%  '
Set-TextValue 'D22' '229.51'
Set-TextValue 'E22' '  -0.30%  '
Set-TextValue 'E23' '  -0.15%  '
Set-TextValue 'D24' '2.44'
Set-TextValue 'E24' '  -0.50%  '
Set-TextValue 'E25' '  +2.85%  '
Set-TextValue 'D26' '9.37'
Set-TextValue 'E26' '  +0.77%  '
Set-TextValue 'D27' '0.142'
Set-TextValue 'E27' '  +14.02%  '
Set-TextValue 'D28' '160.76'
Set-TextValue 'E28' '  -1.28%  '
Set-TextValue 'D29' '19.46'
Set-TextValue 'E29' '  -0.21%  '
Set-TextValue 'E30' '  +1.16%  '
Set-TextValue 'E31' '  -1.34%  '
Set-TextValue 'E32' '  +0.71%  '
Set-TextValue 'D33' '0.0621'
Set-TextValue 'E33' '  -1.74%  '
Set-TextValue 'D34' '4.30'
Set-TextValue 'E34' '  +0.18%  '
Set-TextValue 'B35' 'THORChain'
Set-TextValue 'C35' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D35' '6.12'
Set-TextValue 'E35' '  -2.02%  '
Set-TextValue 'B36' 'BinanceUSD'
Set-TextValue 'C36' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-TextValue 'D36' '1.00'
Set-TextValue 'E36' '  -0.07%  '
Set-TextValue 'D37' '2.27'
Set-TextValue 'E37' '  +5.38%  '
Set-TextValue 'D38' '3.42'
Set-TextValue 'E38' '  +10.73%  '
Set-TextValue 'E39' '  -0.47%  '
Set-TextValue 'D40' '0.101'
Set-TextValue 'E40' '  +3.40%  '
Set-TextValue 'E41' '  -2.33%  '
Set-TextValue 'D42' '0.0212'
Set-TextValue 'E42' '  +1.52%  '
Set-TextValue 'E43' '  -0.95%  '
Set-TextValue 'D44' '16.13'
Set-TextValue 'E44' '  +0.18%  '
Set-TextValue 'D45' '1.362.21'
Set-TextValue 'E45' '  +0.45%  '
Set-TextValue 'E46' '  +0.71%  '
Set-TextValue 'D47' '88.09'
Set-TextValue 'E47' '  +0.22%  '
Set-TextValue 'E48' '  -0.50%  '
Set-TextValue 'D49' '2.84'
Set-TextValue 'E49' '  +0.83%  '
Set-TextValue 'D50' '2.146.63'
Set-TextValue 'E50' '  +1.19%  '
Set-TextValue 'D51' '43.87'
Set-TextValue 'E51' '  -3.35%  '

Write-Output "done"
